# Supplementary tables update: herbicide trade names -> active-ingredient
# names. "Aatrex" becomes "Atrazine" for its very first occurrence in the
# document (the summary table) and "Atrazine-Mesotrione" for every other
# occurrence (the per-site/per-year detail tables). "Clarity" always
# becomes "Dicamba" and "Roundup Powermax" always becomes "Glyphosate".

$d = $word.ActiveDocument
$docEnd = $d.Content.End

function Find-AllPositions($searchText) {
    $positions = @()
    $searchStart = 0
    while ($searchStart -lt $docEnd) {
        $r = $d.Range($searchStart, $docEnd)
        $found = $r.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        $positions += ,@($r.Start, $r.End)
        $searchStart = $r.End
    }
    return $positions
}

# --- Aatrex: first hit -> Atrazine, all later hits -> Atrazine-Mesotrione ---
$aatrexHits = Find-AllPositions("Aatrex")
for ($i = $aatrexHits.Count - 1; $i -ge 0; $i--) {
    $pos = $aatrexHits[$i]
    $r = $d.Range($pos[0], $pos[1])
    if ($i -eq 0) {
        $r.Text = "Atrazine"
    } else {
        $r.Text = "Atrazine-Mesotrione"
    }
}

# --- Clarity -> Dicamba (every occurrence) ---
$clarityHits = Find-AllPositions("Clarity")
for ($i = $clarityHits.Count - 1; $i -ge 0; $i--) {
    $pos = $clarityHits[$i]
    $r = $d.Range($pos[0], $pos[1])
    $r.Text = "Dicamba"
}

# --- Roundup Powermax -> Glyphosate (every occurrence) ---
$roundupHits = Find-AllPositions("Roundup Powermax")
for ($i = $roundupHits.Count - 1; $i -ge 0; $i--) {
    $pos = $roundupHits[$i]
    $r = $d.Range($pos[0], $pos[1])
    $r.Text = "Glyphosate"
}
